$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the header style (A1:H1) back to default (removes bold/border/center style)
$ws.Range("A1:H1").ClearFormats()

# Clear wrap-text style + row height override on H2:H19 / rows 2:19
$ws.Range("H2:H19").ClearFormats()
$ws.Range("2:19").RowHeight = 15

# Update selection to C9
$ws.Range("C9").Select()
